$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8836420178413391
$ws.Range("B1").Value = 1.946284294128418
$ws.Range("C1").Value = 2.793691396713257
$ws.Range("D1").Value = 2.673215389251709
$ws.Range("E1").Value = 0.8693271279335022
